# Update cryptocurrency price/volume snapshot data (GitHub Actions refresh).
# Each targeted cell is forced to text ("@") before assignment so that numeric-
# looking strings (prices, volume codes) are preserved as text, matching the
# original inlineStr cell type, then the number format is restored to General.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.11"
$ws.Range("D2").NumberFormat = "General"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.83"
$ws.Range("D3").NumberFormat = "General"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.333"
$ws.Range("D4").NumberFormat = "General"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.478"
$ws.Range("D6").NumberFormat = "General"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.341"
$ws.Range("D7").NumberFormat = "General"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8103"
$ws.Range("D8").NumberFormat = "General"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8894"
$ws.Range("D9").NumberFormat = "General"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1392"
$ws.Range("D10").NumberFormat = "General"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07349"
$ws.Range("D11").NumberFormat = "General"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03087"
$ws.Range("D12").NumberFormat = "General"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03061"
$ws.Range("D13").NumberFormat = "General"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09348"
$ws.Range("D14").NumberFormat = "General"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.871"
$ws.Range("D15").NumberFormat = "General"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001537"
$ws.Range("D16").NumberFormat = "General"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04717"
$ws.Range("D17").NumberFormat = "General"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006054"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E18").NumberFormat = "General"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005928"
$ws.Range("D19").NumberFormat = "General"

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "UpBots"
$ws.Range("B21").NumberFormat = "General"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("C21").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.007495"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "20UpBotsUBXTBestin24h"
$ws.Range("E21").NumberFormat = "General"

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "HotbitToken"
$ws.Range("B22").NumberFormat = "General"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("C22").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.004654"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "21HotbitTokenHTB"
$ws.Range("E22").NumberFormat = "General"

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("B23").NumberFormat = "General"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("C23").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.00008802"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("E23").NumberFormat = "General"

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "LEO"
$ws.Range("B24").NumberFormat = "General"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C24").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.584"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("E24").NumberFormat = "General"

# Row 25
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("B25").NumberFormat = "General"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C25").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.144"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("E25").NumberFormat = "General"

# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("B26").NumberFormat = "General"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C26").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3180"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("E26").NumberFormat = "General"

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "ProBitToken"
$ws.Range("B27").NumberFormat = "General"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C27").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1318"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "26ProBitTokenPROB"
$ws.Range("E27").NumberFormat = "General"

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "AAXToken"
$ws.Range("B28").NumberFormat = "General"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.2000"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "27AAXTokenAAB"
$ws.Range("E28").NumberFormat = "General"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03806"
$ws.Range("D40").NumberFormat = "General"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.003208"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
$ws.Range("E41").NumberFormat = "General"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002751"
$ws.Range("D42").NumberFormat = "General"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1057"
$ws.Range("D43").NumberFormat = "General"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007855"
$ws.Range("D44").NumberFormat = "General"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005477"
$ws.Range("D45").NumberFormat = "General"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D46").NumberFormat = "General"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5504"
$ws.Range("D47").NumberFormat = "General"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001835"
$ws.Range("D48").NumberFormat = "General"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("D49").NumberFormat = "General"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("D50").NumberFormat = "General"
